# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 130
$ws1.Range("F3").Value = 2151
$ws1.Range("F4").Value = 33
$ws1.Range("F5").Value = 11272
$ws1.Range("F6").Value = 199
$ws1.Range("F9").Value = 11192
$ws1.Range("F12").Value = 56
$ws1.Range("F13").Value = 1731
$ws1.Range("F14").Value = 5599
$ws1.Range("F16").Value = 3456
$ws1.Range("F17").Value = 173

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 130
$ws4.Range("F3").Value = 2151
$ws4.Range("F5").Value = 33
$ws4.Range("F7").Value = 11272
$ws4.Range("F8").Value = 199
$ws4.Range("F11").Value = 11192
$ws4.Range("F14").Value = 56
$ws4.Range("F15").Value = 1731
$ws4.Range("F16").Value = 5599
$ws4.Range("F18").Value = 3456
$ws4.Range("F19").Value = 173
